# Insert a new weekly price record for Perejil (Terminal Hortofrutícola Agro
# Chillán) at row 98, shifting all subsequent records down by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 98:139 down to 99:140, preserving formatting of the row above.
$ws.Rows("98").Insert()

# Populate the newly inserted row 98 with the new record's data.
$ws.Range("A98").Value = 7
$ws.Range("B98").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C98").Value = "Ñuble"
$ws.Range("D98").Value = 45229
$ws.Range("E98").Value = 16
$ws.Range("F98").Value = 100112044
$ws.Range("G98").Value = "Perejil"
$ws.Range("H98").Value = "Sin especificar"
$ws.Range("I98").Value = "Primera"
$ws.Range("J98").Value = 300
$ws.Range("K98").Value = 1500
$ws.Range("L98").Value = 1500
$ws.Range("M98").Value = 1500
$ws.Range("N98").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O98").Value = "Región de Ñuble"
$ws.Range("P98").Value = 1500
$ws.Range("Q98").Value = 1
$ws.Range("R98").Value = "Hortaliza"
